# This script applies updated evaluation metrics to the
# "models_comparison_report.xlsx" workbook, reflecting a re-run of the
# model comparison/evaluation code (per commit message: "[ADD] last code was addded").
#
# Sheet "summary": accuracy / top_2_accuracy / top_3_accuracy per model (rows 2-5).
# Sheet "per_class": count / precision / recall / f1_score per model/class (rows 2-25).

$wb = $excel.ActiveWorkbook
$wsSummary  = $wb.Worksheets.Item("summary")
$wsPerClass = $wb.Worksheets.Item("per_class")

# Update "summary" sheet (rows 2-5: accuracy / top_2_accuracy / top_3_accuracy)
$wsSummary.Range("B2").Value = 0.7320574162679426
$wsSummary.Range("C2").Value = 0.8755980730056763
$wsSummary.Range("D2").Value = 0.9282296895980835

$wsSummary.Range("B3").Value = 0.2727272727272727
$wsSummary.Range("C3").Value = 0.4114832580089569
$wsSummary.Range("D3").Value = 0.6555023789405823

$wsSummary.Range("B4").Value = 0.1961722488038277
$wsSummary.Range("C4").Value = 0.4162679314613342
$wsSummary.Range("D4").Value = 0.6507176756858826

$wsSummary.Range("B5").Value = 0.1100478468899522
$wsSummary.Range("C5").Value = 0.2488038241863251
$wsSummary.Range("D5").Value = 0.3923445045948029


# Update "per_class" sheet (rows 2-25: count / precision / recall / f1_score)
$wsPerClass.Range("D2").Value = 63
$wsPerClass.Range("F2").Value = 0.5079365079365079
$wsPerClass.Range("G2").Value = 0.6736842105263158

$wsPerClass.Range("D3").Value = 25
$wsPerClass.Range("E3").Value = 0.9375
$wsPerClass.Range("F3").Value = 0.6
$wsPerClass.Range("G3").Value = 0.7317073170731707

$wsPerClass.Range("D4").Value = 26
$wsPerClass.Range("E4").Value = 0.5581395348837209
$wsPerClass.Range("F4").Value = 0.9230769230769231
$wsPerClass.Range("G4").Value = 0.6956521739130435

$wsPerClass.Range("D5").Value = 30
$wsPerClass.Range("E5").Value = 0.9032258064516129
$wsPerClass.Range("F5").Value = 0.9333333333333333
$wsPerClass.Range("G5").Value = 0.9180327868852459

$wsPerClass.Range("D6").Value = 56
$wsPerClass.Range("E6").Value = 0.7076923076923077
$wsPerClass.Range("F6").Value = 0.8214285714285714
$wsPerClass.Range("G6").Value = 0.7603305785123967

$wsPerClass.Range("D7").Value = 9
$wsPerClass.Range("E7").Value = 0.3636363636363636
$wsPerClass.Range("F7").Value = 0.8888888888888888
$wsPerClass.Range("G7").Value = 0.5161290322580645

$wsPerClass.Range("D8").Value = 63

$wsPerClass.Range("D9").Value = 25
$wsPerClass.Range("E9").Value = 0.25
$wsPerClass.Range("F9").Value = 0.2
$wsPerClass.Range("G9").Value = 0.2222222222222222

$wsPerClass.Range("D10").Value = 26

$wsPerClass.Range("D11").Value = 30
$wsPerClass.Range("E11").Value = 0.3548387096774194
$wsPerClass.Range("F11").Value = 0.3666666666666666
$wsPerClass.Range("G11").Value = 0.360655737704918

$wsPerClass.Range("D12").Value = 56
$wsPerClass.Range("E12").Value = 0.2697368421052632
$wsPerClass.Range("F12").Value = 0.7321428571428571
$wsPerClass.Range("G12").Value = 0.3942307692307692

$wsPerClass.Range("D13").Value = 9

$wsPerClass.Range("D14").Value = 63
$wsPerClass.Range("E14").Value = 0.2758620689655172
$wsPerClass.Range("F14").Value = 0.253968253968254
$wsPerClass.Range("G14").Value = 0.2644628099173554

$wsPerClass.Range("D15").Value = 25
$wsPerClass.Range("E15").Value = 0.125
$wsPerClass.Range("F15").Value = 0.12
$wsPerClass.Range("G15").Value = 0.1224489795918367

$wsPerClass.Range("D16").Value = 26
$wsPerClass.Range("E16").Value = 0.25
$wsPerClass.Range("F16").Value = 0.03846153846153846
$wsPerClass.Range("G16").Value = 0.06666666666666667

$wsPerClass.Range("D17").Value = 30
$wsPerClass.Range("E17").Value = 0.1707317073170732
$wsPerClass.Range("F17").Value = 0.4666666666666667
$wsPerClass.Range("G17").Value = 0.25

$wsPerClass.Range("D18").Value = 56
$wsPerClass.Range("E18").Value = 0.1842105263157895
$wsPerClass.Range("F18").Value = 0.125
$wsPerClass.Range("G18").Value = 0.148936170212766

$wsPerClass.Range("D19").Value = 9

$wsPerClass.Range("D20").Value = 63

$wsPerClass.Range("D21").Value = 25
$wsPerClass.Range("E21").Value = 0.1363636363636364
$wsPerClass.Range("F21").Value = 0.36
$wsPerClass.Range("G21").Value = 0.1978021978021978

$wsPerClass.Range("D22").Value = 26
$wsPerClass.Range("E22").Value = 0.06451612903225806
$wsPerClass.Range("F22").Value = 0.2307692307692308
$wsPerClass.Range("G22").Value = 0.1008403361344538

$wsPerClass.Range("D23").Value = 30
$wsPerClass.Range("E23").Value = 0.1063829787234043
$wsPerClass.Range("F23").Value = 0.1666666666666667
$wsPerClass.Range("G23").Value = 0.1298701298701299

$wsPerClass.Range("D24").Value = 56
$wsPerClass.Range("E24").Value = 1
$wsPerClass.Range("F24").Value = 0.05357142857142857
$wsPerClass.Range("G24").Value = 0.1016949152542373

$wsPerClass.Range("D25").Value = 9


Write-Host "Updated 'summary' and 'per_class' sheets with refreshed evaluation metrics."
